$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.449.65"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "1.966.69"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.20%  "
$ws.Range("E6").Value = "  -4.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.64"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -9.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -5.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.831"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -10.44%  "
$ws.Range("D15").Value = "2.255.64"
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.85%  "
$ws.Range("E17").Value = "  -5.25%  "
$ws.Range("D18").Value = "1.958.27"
$ws.Range("E18").Value = "  -5.87%  "
$ws.Range("D19").Value = "36.338.66"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("E23").Value = "  -6.99%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("E33").Value = "  -7.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0640"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.47%  "
$ws.Range("E39").Value = "  -15.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0972"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.23%  "
$ws.Range("E44").Value = "  -4.10%  "
$ws.Range("E45").Value = "  -9.83%  "
$ws.Range("E46").Value = "  -8.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "1.348.93"
$ws.Range("E48").Value = "  -3.84%  "
$ws.Range("E49").Value = "  -9.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.40%  "
